$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.71435546875
$ws.Range("C3").Value = 0.71533203125
$ws.Range("C4").Value = 0.72021484375
$ws.Range("C6").Value = 0.7138671875
$ws.Range("C7").Value = 0.712890625
$ws.Range("C8").Value = 0.71435546875
$ws.Range("C9").Value = 0.71533203125
